$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Misc")
$ws.Activate()

$ws.Rows.Item(339).Insert()
$ws.Cells.Item(339,1).Value = "Mickey Moniak 2017 Topps Heritage"
$ws.Cells.Item(339,2).Value = "https://www.topps.com/wp/wp-content/uploads/2017/02/Mickey-Moniak.png"

$excel.ActiveWindow.ScrollRow = 321
$ws.Range("B339").Select()
